# user register process, hand successfully, waiting for verify code finishend
#
# Changes:
#  - wxadmin sheet (sheet2): replace the wxcount/smscount/emailcount header
#    columns (H5:J5) with userlimit/msglimit/memo, widen column H, and move
#    the sheet's remembered selection to C5.
#  - msgcount sheet (sheet3) no longer is the tab-selected sheet (the new
#    wxuser sheet becomes active instead) - handled implicitly below.
#  - a brand new "wxuser" sheet is appended after "msgcount" documenting the
#    wechat-user registration table + the step-by-step registration flow,
#    and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. wxadmin (2nd sheet) - swap out the old count columns for the new
#    userlimit / msglimit / memo columns, resize column H, reposition the
#    saved selection.
# ---------------------------------------------------------------------
$wxadmin = $wb.Worksheets.Item(2)
$wxadmin.Range("H5").Value = "userlimit"
$wxadmin.Range("I5").Value = "msglimit"
$wxadmin.Range("J5").Value = "memo"
$wxadmin.Columns.Item(8).ColumnWidth = 15

$wxadmin.Activate() | Out-Null
$wxadmin.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the new "wxuser" sheet after the last existing sheet (msgcount).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wxuser = $wb.Worksheets.Add($null, $lastSheet)
$wxuser.Name = "wxuser"

$wxuser.Columns.Item(3).ColumnWidth = 11.8

# header row
$wxuser.Range("B4").Value = "id"
$wxuser.Range("C4").Value = "companyid"
$wxuser.Range("D4").Value = "wxid"
$wxuser.Range("E4").Value = "name"
$wxuser.Range("F4").Value = "regcode"
$wxuser.Range("G4").Value = "phone"
$wxuser.Range("H4").Value = "email"
$wxuser.Range("I4").Value = "regtime"
$wxuser.Range("J4").Value = "active"
$wxuser.Range("K4").Value = "step"
$wxuser.Range("L4").Value = "memo"

$wxuser.Range("B5").Value = "主键"
$wxuser.Range("B6").Value = "自增"

$wxuser.Range("B12").Value = "3. 查询是否这个wxid已经注册过，在当前这个wxuser表中是否有这个用户和公司对应的id，并且是active=1，step=9的，如果有有就是注册过，提示已注册，并结束"
$wxuser.Range("C13").Value = "条件是wxid存在，且active=1"

$wxuser.Range("B16").Value = "4. 如果查出来active等于1，并且step=1，则该用户在等待注册码，跳转至输入注册码页面"
$wxuser.Range("B17").Value = "5.  如果这个用户未注册，在wxadmin中查询这个公司id及购买的userlimit的总数，没有超过poinfo的规定限制，跳转至用户输入信息的页面，点击提交把用户填写的信息保存在表里面，并把注册码保存在本表中，把step设置成1"

$wxuser.Range("B23").Value = "1. 首先从wxadmin表中查询出来该公司id对应的管理员的wxid（微信账号）"
$wxuser.Range("B24").Value = "2. 如果没有这个公司id，或者enable未0，提示用户管理员还未注册，该公司账号还没有生效，并结束"

$wxuser.Range("B28").Value = "6. 跳转至等待输入注册码的页面"
$wxuser.Range("B29").Value = "7. 给管理员的微信发送一个模板消息，模板消息中带有随机生成的注册码"
$wxuser.Range("B30").Value = "8. 请这个微信用户找管理员要这个注册码，如果输入正确，则账号激活，step设置成9"

# ---------------------------------------------------------------------
# 3. Make the new sheet the active / tab-selected sheet with its saved
#    selection on C14.
# ---------------------------------------------------------------------
$wxuser.Activate() | Out-Null
$wxuser.Range("C14").Select() | Out-Null
